$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current data (rows 3-14) for each algorithm name, keyed by name,
# so we can re-write it in the new order without losing any values.
$data = @{}
for ($r = 3; $r -le 14; $r++) {
    $name = $ws.Range("A$r").Value2
    $data[$name] = @(
        $ws.Range("B$r").Value2,
        $ws.Range("C$r").Value2,
        $ws.Range("D$r").Value2,
        $ws.Range("E$r").Value2,
        $ws.Range("F$r").Value2
    )
}

# New row order (rows 3-14) after the re-sort.
$newOrder = @(
    "cem",
    "wachter",
    "face-epsilon",
    "cchvae",
    "revise",
    "ar",
    "cruds",
    "face-knn",
    "clue",
    "cem-vae",
    "gs",
    "dice"
)

$r = 3
foreach ($name in $newOrder) {
    $vals = $data[$name]
    $ws.Range("A$r").Value = $name
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("F$r").Value = $vals[4]
    $r++
}
